# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) values for the relevant rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    3  = 0
    4  = 1
    5  = 1
    6  = 1
    7  = 2
    8  = 1
    9  = 1
    10 = 1
    11 = 2
    12 = 2
    13 = 1
    14 = 3
    15 = 1
    16 = 1
    17 = 1
    19 = 1
    20 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
